$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Cells.Item(2, 1).Value = '{''Hobby'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(2, 2).Value = 20675
$ws.Cells.Item(2, 3).Value = 9822.513984659707
$ws.Cells.Item(2, 4).Value = -3140.735516262439
$ws.Cells.Item(3, 1).Value = '{''UndergradMajor'': np.int64(2), ''Student'': np.int64(1)}'
$ws.Cells.Item(3, 2).Value = 15077
$ws.Cells.Item(3, 3).Value = 20114.57186124945
$ws.Cells.Item(3, 4).Value = 7151.322360327307
$ws.Cells.Item(4, 1).Value = '{''Student'': np.int64(1), ''DevType'': np.int64(2)}'
$ws.Cells.Item(4, 2).Value = 15397
$ws.Cells.Item(4, 3).Value = 12458.93044835459
$ws.Cells.Item(4, 4).Value = -504.3190525675527
$ws.Cells.Item(5, 1).Value = '{''Gender'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Cells.Item(5, 2).Value = 23188
$ws.Cells.Item(5, 3).Value = 10204.27301983108
$ws.Cells.Item(5, 4).Value = -2758.976481091066
$ws.Cells.Item(6, 1).Value = '{''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(6, 2).Value = 23264
$ws.Cells.Item(6, 3).Value = 10145.98630245002
$ws.Cells.Item(6, 4).Value = -2817.263198472123
$ws.Cells.Item(7, 1).Value = '{''Dependents'': np.int64(2), ''Student'': np.int64(1)}'
$ws.Cells.Item(7, 2).Value = 16812
$ws.Cells.Item(7, 3).Value = 14915.38565374814
$ws.Cells.Item(7, 4).Value = 1952.136152825999
$ws.Cells.Item(8, 1).Value = '{''HDI'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Cells.Item(8, 2).Value = 21710
$ws.Cells.Item(8, 3).Value = 10343.01875536946
$ws.Cells.Item(8, 4).Value = -2620.23074555269
$ws.Cells.Item(9, 1).Value = '{''Gender'': np.int64(1), ''UndergradMajor'': np.int64(2)}'
$ws.Cells.Item(9, 2).Value = 17057
$ws.Cells.Item(9, 3).Value = 20268.49410164942
$ws.Cells.Item(9, 4).Value = 7305.244600727277
$ws.Cells.Item(10, 1).Value = '{''UndergradMajor'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(10, 2).Value = 16908
$ws.Cells.Item(10, 3).Value = 19053.63946379748
$ws.Cells.Item(10, 4).Value = 6090.389962875339
$ws.Cells.Item(11, 1).Value = '{''UndergradMajor'': np.int64(2), ''HDI'': np.int64(1)}'
$ws.Cells.Item(11, 2).Value = 15406
$ws.Cells.Item(11, 3).Value = 19129.56707097096
$ws.Cells.Item(11, 4).Value = 6166.317570048814
$ws.Cells.Item(12, 1).Value = '{''Gender'': np.int64(1), ''DevType'': np.int64(2)}'
$ws.Cells.Item(12, 2).Value = 16852
$ws.Cells.Item(12, 3).Value = 11505.16324489906
$ws.Cells.Item(12, 4).Value = -1458.086256023089
$ws.Cells.Item(13, 1).Value = '{''DevType'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(13, 2).Value = 16631
$ws.Cells.Item(13, 3).Value = 10372.18434513095
$ws.Cells.Item(13, 4).Value = -2591.065155791193
$ws.Cells.Item(14, 1).Value = '{''HDI'': np.int64(1), ''DevType'': np.int64(2)}'
$ws.Cells.Item(14, 2).Value = 15175
$ws.Cells.Item(14, 3).Value = 11257.39659382552
$ws.Cells.Item(14, 4).Value = -1705.852907096631
$ws.Cells.Item(15, 1).Value = '{''Gender'': np.int64(1), ''HoursComputer'': np.int64(2)}'
$ws.Cells.Item(15, 2).Value = 15415
$ws.Cells.Item(15, 3).Value = 11542.48855671203
$ws.Cells.Item(15, 4).Value = -1420.760944210115
$ws.Cells.Item(16, 1).Value = '{''HoursComputer'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(16, 2).Value = 15402
$ws.Cells.Item(16, 3).Value = 10116.13933820374
$ws.Cells.Item(16, 4).Value = -2847.110162718405
$ws.Cells.Item(17, 1).Value = '{''Hobby'': np.int64(1)}'
$ws.Cells.Item(17, 2).Value = 24271
$ws.Cells.Item(17, 3).Value = 11442.32763080001
$ws.Cells.Item(17, 4).Value = -1520.921870122134
$ws.Cells.Item(18, 1).Value = '{''Student'': np.int64(1)}'
$ws.Cells.Item(18, 2).Value = 25010
$ws.Cells.Item(18, 3).Value = 12118.41632988031
$ws.Cells.Item(18, 4).Value = -844.833171041837
$ws.Cells.Item(19, 1).Value = '{''UndergradMajor'': np.int64(2)}'
$ws.Cells.Item(19, 2).Value = 18141
$ws.Cells.Item(19, 3).Value = 20428.79458652144
$ws.Cells.Item(19, 4).Value = 7465.545085599289
$ws.Cells.Item(20, 1).Value = '{''DevType'': np.int64(2)}'
$ws.Cells.Item(20, 2).Value = 17900
$ws.Cells.Item(20, 3).Value = 12666.05425130746
$ws.Cells.Item(20, 4).Value = -297.1952496146841
$ws.Cells.Item(21, 1).Value = '{''HoursComputer'': np.int64(2)}'
$ws.Cells.Item(21, 2).Value = 16618
$ws.Cells.Item(21, 3).Value = 11665.80607096392
$ws.Cells.Item(21, 4).Value = -1297.443429958223
$ws.Cells.Item(22, 1).Value = '{''Gender'': np.int64(1)}'
$ws.Cells.Item(22, 2).Value = 27355
$ws.Cells.Item(22, 3).Value = 11218.84231805306
$ws.Cells.Item(22, 4).Value = -1744.407182869083
$ws.Cells.Item(23, 1).Value = '{''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(23, 2).Value = 27379
$ws.Cells.Item(23, 3).Value = 11090.60929182074
$ws.Cells.Item(23, 4).Value = -1872.640209101404
$ws.Cells.Item(24, 1).Value = '{''Dependents'': np.int64(2)}'
$ws.Cells.Item(24, 2).Value = 20567
$ws.Cells.Item(24, 3).Value = 15977.15447741413
$ws.Cells.Item(24, 4).Value = 3013.904976491984
$ws.Cells.Item(25, 1).Value = '{''Age'': np.int64(3)}'
$ws.Cells.Item(25, 2).Value = 15063
$ws.Cells.Item(25, 3).Value = 13540.92143844457
$ws.Cells.Item(25, 4).Value = 577.6719375224238
$ws.Cells.Item(26, 1).Value = '{''HDI'': np.int64(1)}'
$ws.Cells.Item(26, 2).Value = 25023
$ws.Cells.Item(26, 3).Value = 11250.47982698071
$ws.Cells.Item(26, 4).Value = -1712.76967394144
$ws.Cells.Item(27, 1).Value = '{''Hobby'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Cells.Item(27, 2).Value = 20375
$ws.Cells.Item(27, 3).Value = 10196.88929291665
$ws.Cells.Item(27, 4).Value = -2766.360208005492
$ws.Cells.Item(28, 1).Value = '{''UndergradMajor'': np.int64(2), ''Hobby'': np.int64(1)}'
$ws.Cells.Item(28, 2).Value = 15058
$ws.Cells.Item(28, 3).Value = 18331.65921398864
$ws.Cells.Item(28, 4).Value = 5368.409713066496
$ws.Cells.Item(29, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1)}'
$ws.Cells.Item(29, 2).Value = 22770
$ws.Cells.Item(29, 3).Value = 9793.201298325481
$ws.Cells.Item(29, 4).Value = -3170.048202596665
$ws.Cells.Item(30, 1).Value = '{''Hobby'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(30, 2).Value = 22529
$ws.Cells.Item(30, 3).Value = 8905.820842777424
$ws.Cells.Item(30, 4).Value = -4057.428658144721
$ws.Cells.Item(31, 1).Value = '{''Dependents'': np.int64(2), ''Hobby'': np.int64(1)}'
$ws.Cells.Item(31, 2).Value = 17063
$ws.Cells.Item(31, 3).Value = 13434.95249482262
$ws.Cells.Item(31, 4).Value = 471.7029939004733
$ws.Cells.Item(32, 1).Value = '{''Gender'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(32, 2).Value = 25910
$ws.Cells.Item(32, 3).Value = 10263.12279291961
$ws.Cells.Item(32, 4).Value = -2700.126708002535
$ws.Cells.Item(33, 1).Value = '{''Gender'': np.int64(1), ''Dependents'': np.int64(2)}'
$ws.Cells.Item(33, 2).Value = 18813
$ws.Cells.Item(33, 3).Value = 14017.15883098588
$ws.Cells.Item(33, 4).Value = 1053.909330063734
$ws.Cells.Item(34, 1).Value = '{''Gender'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(34, 2).Value = 23143
$ws.Cells.Item(34, 3).Value = 9283.890006454172
$ws.Cells.Item(34, 4).Value = -3679.359494467973
$ws.Cells.Item(35, 1).Value = '{''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(35, 2).Value = 18715
$ws.Cells.Item(35, 3).Value = 13935.86347048962
$ws.Cells.Item(35, 4).Value = 972.6139695674774
$ws.Cells.Item(36, 1).Value = '{''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(36, 2).Value = 23146
$ws.Cells.Item(36, 3).Value = 9368.314941529816
$ws.Cells.Item(36, 4).Value = -3594.93455939233
$ws.Cells.Item(37, 1).Value = '{''Dependents'': np.int64(2), ''HDI'': np.int64(1)}'
$ws.Cells.Item(37, 2).Value = 17336
$ws.Cells.Item(37, 3).Value = 14861.36267076715
$ws.Cells.Item(37, 4).Value = 1898.113169845003
$ws.Cells.Item(38, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Cells.Item(38, 2).Value = 19121
$ws.Cells.Item(38, 3).Value = 8516.623745728855
$ws.Cells.Item(38, 4).Value = -4446.62575519329
$ws.Cells.Item(39, 1).Value = '{''Hobby'': np.int64(1), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(39, 2).Value = 18956
$ws.Cells.Item(39, 3).Value = 7726.243805460593
$ws.Cells.Item(39, 4).Value = -5237.005695461553
$ws.Cells.Item(40, 1).Value = '{''Hobby'': np.int64(1), ''Student'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(40, 2).Value = 17771
$ws.Cells.Item(40, 3).Value = 8530.523205969457
$ws.Cells.Item(40, 4).Value = -4432.726294952688
$ws.Cells.Item(41, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(41, 2).Value = 21544
$ws.Cells.Item(41, 3).Value = 8228.149376149711
$ws.Cells.Item(41, 4).Value = -4735.100124772434
$ws.Cells.Item(42, 1).Value = '{''Gender'': np.int64(1), ''Dependents'': np.int64(2), ''Hobby'': np.int64(1)}'
$ws.Cells.Item(42, 2).Value = 15834
$ws.Cells.Item(42, 3).Value = 11719.31701495631
$ws.Cells.Item(42, 4).Value = -1243.932485965832
$ws.Cells.Item(43, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(43, 2).Value = 19359
$ws.Cells.Item(43, 3).Value = 7940.548127036916
$ws.Cells.Item(43, 4).Value = -5022.701373885229
$ws.Cells.Item(44, 1).Value = '{''Dependents'': np.int64(2), ''Hobby'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(44, 2).Value = 15532
$ws.Cells.Item(44, 3).Value = 10474.02650894398
$ws.Cells.Item(44, 4).Value = -2489.222991978169
$ws.Cells.Item(45, 1).Value = '{''Hobby'': np.int64(1), ''SexualOrientation'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(45, 2).Value = 19133
$ws.Cells.Item(45, 3).Value = 7281.607171741128
$ws.Cells.Item(45, 4).Value = -5681.642329181018
$ws.Cells.Item(46, 1).Value = '{''Gender'': np.int64(1), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(46, 2).Value = 22014
$ws.Cells.Item(46, 3).Value = 9367.700301271154
$ws.Cells.Item(46, 4).Value = -3595.549199650992
$ws.Cells.Item(47, 1).Value = '{''Hobby'': np.int64(1), ''SexualOrientation'': np.int64(1), ''Gender'': np.int64(1), ''HDI'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Cells.Item(47, 2).Value = 15779
$ws.Cells.Item(47, 3).Value = 5234.178807957504
$ws.Cells.Item(47, 4).Value = -7729.070692964642
$ws.Cells.Item(48, 1).Value = '{''Gender'': np.int64(1), ''Dependents'': np.int64(2), ''Student'': np.int64(1)}'
$ws.Cells.Item(48, 2).Value = 15339
$ws.Cells.Item(48, 3).Value = 12445.08873558062
$ws.Cells.Item(48, 4).Value = -518.1607653415285
$ws.Cells.Item(49, 1).Value = '{''Gender'': np.int64(1), ''HDI'': np.int64(1), ''Student'': np.int64(1)}'
$ws.Cells.Item(49, 2).Value = 20079
$ws.Cells.Item(49, 3).Value = 8320.806949419126
$ws.Cells.Item(49, 4).Value = -4642.44255150302
$ws.Cells.Item(50, 1).Value = '{''Dependents'': np.int64(2), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(50, 2).Value = 15301
$ws.Cells.Item(50, 3).Value = 12493.22924053899
$ws.Cells.Item(50, 4).Value = -470.0202603831513
$ws.Cells.Item(51, 1).Value = '{''HDI'': np.int64(1), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(51, 2).Value = 20128
$ws.Cells.Item(51, 3).Value = 8058.492947675052
$ws.Cells.Item(51, 4).Value = -4904.756553247094
$ws.Cells.Item(52, 1).Value = '{''Gender'': np.int64(1), ''UndergradMajor'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(52, 2).Value = 16152
$ws.Cells.Item(52, 3).Value = 19730.73379293201
$ws.Cells.Item(52, 4).Value = 6767.48429200986
$ws.Cells.Item(53, 1).Value = '{''Gender'': np.int64(1), ''DevType'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(53, 2).Value = 15929
$ws.Cells.Item(53, 3).Value = 10271.32109283204
$ws.Cells.Item(53, 4).Value = -2691.928408090102
$ws.Cells.Item(54, 1).Value = '{''Gender'': np.int64(1), ''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(54, 2).Value = 17563
$ws.Cells.Item(54, 3).Value = 12896.98458729801
$ws.Cells.Item(54, 4).Value = -66.26491362413253
$ws.Cells.Item(55, 1).Value = '{''Gender'': np.int64(1), ''HDI'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(55, 2).Value = 21900
$ws.Cells.Item(55, 3).Value = 8524.583644570461
$ws.Cells.Item(55, 4).Value = -4438.665856351685
$ws.Cells.Item(56, 1).Value = '{''Gender'': np.int64(1), ''Dependents'': np.int64(2), ''HDI'': np.int64(1)}'
$ws.Cells.Item(56, 2).Value = 15786
$ws.Cells.Item(56, 3).Value = 12632.85656038428
$ws.Cells.Item(56, 4).Value = -330.392940537864
$ws.Cells.Item(57, 1).Value = '{''Dependents'': np.int64(2), ''SexualOrientation'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(57, 2).Value = 15673
$ws.Cells.Item(57, 3).Value = 12988.61598577243
$ws.Cells.Item(57, 4).Value = 25.36648485028491
$ws.Cells.Item(58, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(58, 2).Value = 18136
$ws.Cells.Item(58, 3).Value = 7283.726461989329
$ws.Cells.Item(58, 4).Value = -5679.523038932816
$ws.Cells.Item(59, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''Student'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(59, 2).Value = 16648
$ws.Cells.Item(59, 3).Value = 6709.311194402443
$ws.Cells.Item(59, 4).Value = -6253.938306519703
$ws.Cells.Item(60, 1).Value = '{''Hobby'': np.int64(1), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(60, 2).Value = 16486
$ws.Cells.Item(60, 3).Value = 5663.715061677382
$ws.Cells.Item(60, 4).Value = -7299.534439244764
$ws.Cells.Item(61, 1).Value = '{''Gender'': np.int64(1), ''Hobby'': np.int64(1), ''SexualOrientation'': np.int64(1), ''HDI'': np.int64(1)}'
$ws.Cells.Item(61, 2).Value = 18306
$ws.Cells.Item(61, 3).Value = 6562.835401713361
$ws.Cells.Item(61, 4).Value = -6400.414099208785
$ws.Cells.Item(62, 1).Value = '{''Gender'': np.int64(1), ''HDI'': np.int64(1), ''Student'': np.int64(1), ''SexualOrientation'': np.int64(1)}'
$ws.Cells.Item(62, 2).Value = 19043
$ws.Cells.Item(62, 3).Value = 7319.305231565133
$ws.Cells.Item(62, 4).Value = -5643.944269357013
